$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear G6 (was "گیت و گیت هاب (2)") - content moved down to G10
$ws.Range("G6").Clear()

# Set G7 with new workshop text
$ws.Range("G7").Value = "نحوه ذخیره سازی مقادیر در کامپیوتر (1)"

# Set G10 with text that used to be in G6
$ws.Range("G10").Value = "گیت و گیت هاب (2)"

# Clear G13 (was "واسط گرافیکی (1)") - content moved down to G20
$ws.Range("G13").Clear()

# G20 now holds what used to be in G13 ("واسط گرافیکی (1)"), replacing "داکر (1)"
$ws.Range("G20").Value = "واسط گرافیکی (1)"

# Column G grew wider to fit the new longer text
$ws.Columns("G:G").AutoFit()

# Update the view: zoom to 110%, reset scroll position, select G10
$ws.Application.ActiveWindow.Zoom = 110
$ws.Range("G10").Select()
